# Update countries & provincias Spain
# Applies the daily COVID data refresh: updates the "last updated" timestamp,
# refreshes case counters for a number of countries, and re-ranks rows whose
# "Casos totales" changed enough to swap their position in the (descending
# sorted) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" footer text (last row of the sheet) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 1 de Octubre de 2020 a las 02:22"

# --- 2. Row data refresh / re-rank. Each entry is:
#        (row, Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes)
$rows = @(
    @(4,   "Estados Unidos",       7446825, 40472, 4688859, 2546253, 0, 928, 211713),
    @(6,   "Brasil",               4813586, 33269, 4180376,  489248, 0, 952, 143962),
    @(29,  "Canada",                158758,  1797,  134971,   14490, 0,   6,   9297),
    @(37,  "Panama",                112595,   742,   89061,   21162, 0,   8,   2372),
    @(38,  "Republica Dominicana",  112209,   309,   87629,   22475, 0,   4,   2105),
    @(102, "Guinea",                 10652,    18,    9996,     590, 0,   0,     66),
    @(125, "Congo",                   5089,    81,    3887,    1113, 0,   0,     89),
    @(126, "Hong Kong",               5088,     8,    4827,     156, 0,   0,    105),
    @(127, "Guinea Ecuatorial",       5030,     0,    4769,     178, 0,   0,     83),
    @(129, "Surinam",                 4877,    14,    4695,      78, 0,   0,    104),
    @(158, "Polinesia Francesa",      1852,   124,    1504,     341, 0,   0,      7),
    @(159, "Nueva Zelanda",           1836,     1,    1767,      44, 0,   0,     25),
    @(160, "Letonia",                 1824,    95,    1307,     480, 0,   0,     37),
    @(161, "Togo",                    1784,    25,    1348,     388, 0,   0,     48),
    @(162, "Republica de Chipre",     1755,    12,    1369,     364, 0,   0,     22),
    @(183, "Mauricio",                 381,     0,     344,      27, 0,   0,     10),
    @(184, "Eritrea",                  375,     0,     341,      34, 0,   0,      0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
